# Update dates and links (Session1.pptx)
#
# Slide 1: "Hilary 2023" -> "Trinity 2023"
# Slide 2: GitHub URL run merges with the old "/MSD_R_course_HT2023" run,
#          text becomes "https://github.com/sraorao/MSD_R_course_TT2023" and
#          the now-empty trailing run is removed; "Screen " + "sharing" runs
#          merge into a single "Screen sharing" run.
#
# NOTE: this COM-interop engine implements TextRange.Text assignment as a
# longest-common-*suffix* splice (it keeps any trailing substring shared by
# the old and new text as a separate, unedited run). To land a clean,
# single-run result that matches how a human retyping the text would end up,
# we first push the target sub-range through a throwaway value that shares
# no suffix with either the old or the final text, then set the final text
# in a second pass. That avoids leaving spurious extra <a:r> runs behind.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 1 - "Hilary 2023" -> "Trinity 2023"
# ---------------------------------------------------------------------
$s1 = $p.Slides.Item(1)
$shp1 = $s1.Shapes.Item(2)
$tr1 = $shp1.TextFrame.TextRange
$datePara = $tr1.Paragraphs(2, 1)

$datePara.Text = "~~~~~~~~~~~~"
$datePara = $tr1.Paragraphs(2, 1)
$datePara.Text = "Trinity 2023"

# ---------------------------------------------------------------------
# Slide 2
# ---------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$shp2 = $s2.Shapes.Item(1)
$tr2 = $shp2.TextFrame.TextRange

# --- Paragraph 2: GitHub repository link ---
# Current runs: "GitHub repository " | "https://github.com/sraorao" | "/MSD_R_course_HT2023"
$para2 = $tr2.Paragraphs(2, 1)

# Drop the trailing "/MSD_R_course_HT2023" run entirely.
$oldSuffixRun = $para2.Characters(19 + 26, 20)
$oldSuffixRun.Text = ""

# Grow the URL run's text in place (keeps its rPr / hyperlink) into the
# full new address.
$para2 = $tr2.Paragraphs(2, 1)
$urlRun = $para2.Characters(19, 26)
$urlRun.Text = "zzzzzzzzzzzzzzzzzzzzzzzzzz"
$para2 = $tr2.Paragraphs(2, 1)
$urlRun = $para2.Characters(19, 26)
$urlRun.Text = "https://github.com/sraorao/MSD_R_course_TT2023"

# --- Paragraph 3: "Screen " + "sharing" -> "Screen sharing" ---
$para3 = $tr2.Paragraphs(3, 1)

# Remove the "Screen " run.
$screenRun = $para3.Characters(1, 7)
$screenRun.Text = ""

# Expand the "sharing" run's text (keeps its rPr with dirty="0").
$para3 = $tr2.Paragraphs(3, 1)
$sharingRun = $para3.Characters(1, 7)
$sharingRun.Text = "yyyyyyy"
$para3 = $tr2.Paragraphs(3, 1)
$sharingRun = $para3.Characters(1, 7)
$sharingRun.Text = "Screen sharing"
